$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-06 Thursday", "2025-02-07 Friday"),
    @("42×29=1218", "85×22=1870"),
    @("41×14=574", "50×67=3350"),
    @("20×90=1800", "52×18=936"),
    @("11×14=154", "72×59=4248"),
    @("47×62=2914", "21×97=2037"),
    @("78×44=3432", "75×67=5025"),
    @("49×39=1911", "22×29=638"),
    @("18×76=1368", "33×43=1419"),
    @("30×95=2850", "29×35=1015"),
    @("60×56=3360", "47×86=4042"),
    @("27×66=1782", "75×15=1125"),
    @("40×43=1720", "18×60=1080"),
    @("40×25=1000", "54×55=2970"),
    @("84×45=3780", "42×34=1428"),
    @("81×65=5265", "93×13=1209"),
    @("55×74=4070", "27×90=2430"),
    @("79×96=7584", "49×79=3871"),
    @("13×50=650", "57×72=4104"),
    @("71×36=2556", "47×17=799"),
    @("41×15=615", "45×68=3060"),
    @("22×32=704", "92×95=8740"),
    @("63×79=4977", "86×74=6364"),
    @("79×24=1896", "32×47=1504"),
    @("42×45=1890", "69×63=4347"),
    @("48×82=3936", "69×29=2001")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "done"
